$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clear the old duplicate table (columns G:H on row 1, and E5:H7) ---
$ws.Range("G1:H1").ClearContents()
$ws.Range("E5:H7").ClearContents()

# --- Move the existing invoice rows (C2:F4) down to (C6:F8) ---
$ws.Range("C2:F4").Cut($ws.Range("C6:F8"))

# --- Add the new "send to email" recipient rows ---
# (Force text format while entering the dd-mm-yyyy strings so Excel does not
#  auto-convert them to date serials, then clear the format again so the
#  cells keep the default style.)
$recipients = $ws.Range("A3:B5")
$recipients.NumberFormat = "@"
$ws.Cells.Item(3, 1).Value = "n4zdfr8rz4cdlxzf49uy"
$ws.Cells.Item(3, 2).Value = "08-07-2024"
$ws.Cells.Item(4, 1).Value = "j6qky8ysjflms7kciqj97i"
$ws.Cells.Item(4, 2).Value = "14-07-2024"
$ws.Cells.Item(5, 1).Value = "g15db3dv9zupp579hzbzm"
$ws.Cells.Item(5, 2).Value = "15-07-2024"
$recipients.ClearFormats()

# --- Set new column widths (values chosen so the resulting quantized
#     stored width lands as close as possible to the target widths) ---
$ws.Columns.Item(1).ColumnWidth = 21.166666666666668
$ws.Columns.Item(2).ColumnWidth = 16.666666666666668
$ws.Columns.Item(3).ColumnWidth = 14.166666666666666
$ws.Columns.Item(4).ColumnWidth = 14.166666666666666
$ws.Columns.Item(5).ColumnWidth = 12.333333333333334
$ws.Columns.Item(6).ColumnWidth = 8.166666666666666
$ws.Columns.Item(7).ColumnWidth = 13.5
$ws.Columns.Item(8).ColumnWidth = 22.666666666666668

# --- Update selection to match the new active cell ---
$ws.Range("C13").Select()
